$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows above row 229, shifting the existing rows
# (229-258) down to (232-261).
$ws.Range("A229:T231").EntireRow.Insert()

# The 3 new rows share the same "static" columns (A,B,C,F,G,H,I,J) as the
# rest of the Agricola del Norte S.A. de Arica / Mango block, which is now
# visible starting at row 232 (old row 229, shifted down by 3).
$staticCols = @("A","B","C","F","G","H","I","J")
foreach ($col in $staticCols) {
    $src = $ws.Range($col + "232").Value()
    $ws.Range($col + "229").Value = $src
    $ws.Range($col + "230").Value = $src
    $ws.Range($col + "231").Value = $src
}

# Keep the same date number format on the new rows' "Fecha" column as the
# rest of the block (style index 2 / numFmtId 165).
$ws.Range("D229:D231").NumberFormat = $ws.Range("D232").NumberFormat()

# New row 229 data
$ws.Range("D229").Value = 45034
$ws.Range("E229").Value = 15
$ws.Range("K229").Value = "Piqueño"
$ws.Range("L229").Value = "Primera"
$ws.Range("M229").Value = 100
$ws.Range("N229").Value = 9000
$ws.Range("O229").Value = 10000
$ws.Range("P229").Value = 9500
$ws.Range("Q229").Value = "$/caja 10 kilos"
$ws.Range("R229").Value = "Región de Arica y Parinacota"
$ws.Range("S229").Value = 950
$ws.Range("T229").Value = 10

# New row 230 data
$ws.Range("D230").Value = 45034
$ws.Range("E230").Value = 15
$ws.Range("K230").Value = "Piqueño"
$ws.Range("L230").Value = "Segunda"
$ws.Range("M230").Value = 160
$ws.Range("N230").Value = 7000
$ws.Range("O230").Value = 8000
$ws.Range("P230").Value = 7500
$ws.Range("Q230").Value = "$/caja 10 kilos"
$ws.Range("R230").Value = "Región de Arica y Parinacota"
$ws.Range("S230").Value = 750
$ws.Range("T230").Value = 10

# New row 231 data
$ws.Range("D231").Value = 45034
$ws.Range("E231").Value = 15
$ws.Range("K231").Value = "Piqueño"
$ws.Range("L231").Value = "Tercera"
$ws.Range("M231").Value = 160
$ws.Range("N231").Value = 5000
$ws.Range("O231").Value = 6000
$ws.Range("P231").Value = 5500
$ws.Range("Q231").Value = "$/caja 10 kilos"
$ws.Range("R231").Value = "Región de Arica y Parinacota"
$ws.Range("S231").Value = 550
$ws.Range("T231").Value = 10
